# Update the log write mode: refresh run_time, max_er, and per-iteration
# log values recorded for each simulation run on the active worksheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 3).Value = 0.3260133266448975
$ws.Cells.Item(2, 5).Value = 55.19530194853905
$ws.Cells.Item(2, 6).Value = 0.00193708617319336
$ws.Cells.Item(2, 7).Value = 0.001614282915651162
$ws.Cells.Item(2, 8).Value = 0.001474663765395419
$ws.Cells.Item(2, 9).Value = 0.001474663765395419
$ws.Cells.Item(2, 10).Value = 0.001435007657782635
$ws.Cells.Item(2, 11).Value = 0.001368985721015087
$ws.Cells.Item(2, 12).Value = 0.001342061859284133
$ws.Cells.Item(2, 13).Value = 0.001247027965317697
$ws.Cells.Item(2, 14).Value = 0.001247027965317697
$ws.Cells.Item(2, 15).Value = 0.001161042037927625
$ws.Cells.Item(2, 16).Value = 0.001161042037927625
$ws.Cells.Item(2, 17).Value = 0.001122052931221521
$ws.Cells.Item(2, 18).Value = 0.001097907715602116
$ws.Cells.Item(2, 19).Value = 0.001097907715602116
$ws.Cells.Item(2, 20).Value = 0.001097907715602116
$ws.Cells.Item(2, 21).Value = 0.001097907715602116
$ws.Cells.Item(2, 22).Value = 0.001096096261867095
$ws.Cells.Item(2, 23).Value = 0.001094757477420461
$ws.Cells.Item(2, 24).Value = 0.001080757880222699
$ws.Cells.Item(2, 25).Value = 0.001075931811862359
$ws.Cells.Item(3, 3).Value = 0.2740449905395508
$ws.Cells.Item(3, 5).Value = 51.90189426400138
$ws.Cells.Item(3, 6).Value = 0.00193708617319336
$ws.Cells.Item(3, 7).Value = 0.001588698868452422
$ws.Cells.Item(3, 8).Value = 0.001449352956172938
$ws.Cells.Item(3, 9).Value = 0.001327057453588113
$ws.Cells.Item(3, 10).Value = 0.001327057453588113
$ws.Cells.Item(3, 11).Value = 0.001285011079778936
$ws.Cells.Item(3, 12).Value = 0.001253697775171577
$ws.Cells.Item(3, 13).Value = 0.001246499520321984
$ws.Cells.Item(3, 14).Value = 0.001246499520321984
$ws.Cells.Item(3, 15).Value = 0.001167888378278615
$ws.Cells.Item(3, 16).Value = 0.001144224340544481
$ws.Cells.Item(3, 17).Value = 0.001104392785848893
$ws.Cells.Item(3, 18).Value = 0.001094938547026175
$ws.Cells.Item(3, 19).Value = 0.001073635332183636
$ws.Cells.Item(3, 20).Value = 0.00105479369664574
$ws.Cells.Item(3, 21).Value = 0.001037888920513721
$ws.Cells.Item(3, 22).Value = 0.001037888920513721
$ws.Cells.Item(3, 23).Value = 0.001023560124360822
$ws.Cells.Item(3, 24).Value = 0.001016834911464663
$ws.Cells.Item(3, 25).Value = 0.001011732831656947
$ws.Cells.Item(4, 3).Value = 0.2221107482910156
$ws.Cells.Item(4, 5).Value = 53.1604029801274
$ws.Cells.Item(4, 6).Value = 0.00193708617319336
$ws.Cells.Item(4, 7).Value = 0.001474373738155936
$ws.Cells.Item(4, 8).Value = 0.00142741662757629
$ws.Cells.Item(4, 9).Value = 0.001383611703336198
$ws.Cells.Item(4, 10).Value = 0.001361862594219478
$ws.Cells.Item(4, 11).Value = 0.001327506828734737
$ws.Cells.Item(4, 12).Value = 0.001256940992374843
$ws.Cells.Item(4, 13).Value = 0.001256940992374843
$ws.Cells.Item(4, 14).Value = 0.001256940992374843
$ws.Cells.Item(4, 15).Value = 0.001256940992374843
$ws.Cells.Item(4, 16).Value = 0.001195443298494216
$ws.Cells.Item(4, 17).Value = 0.001195443298494216
$ws.Cells.Item(4, 18).Value = 0.001154150309460297
$ws.Cells.Item(4, 19).Value = 0.001139217272490699
$ws.Cells.Item(4, 20).Value = 0.001098253733138509
$ws.Cells.Item(4, 21).Value = 0.001087842380318435
$ws.Cells.Item(4, 22).Value = 0.001080728982990719
$ws.Cells.Item(4, 23).Value = 0.001057865891180823
$ws.Cells.Item(4, 24).Value = 0.001037905199908826
$ws.Cells.Item(4, 25).Value = 0.001036265165304628
$ws.Cells.Item(5, 3).Value = 0.2858626842498779
$ws.Cells.Item(5, 5).Value = 55.91586791694499
$ws.Cells.Item(5, 6).Value = 0.00193708617319336
$ws.Cells.Item(5, 7).Value = 0.001607069092659935
$ws.Cells.Item(5, 8).Value = 0.001544747418441592
$ws.Cells.Item(5, 9).Value = 0.001397212136701191
$ws.Cells.Item(5, 10).Value = 0.001397212136701191
$ws.Cells.Item(5, 11).Value = 0.001397212136701191
$ws.Cells.Item(5, 12).Value = 0.001334379282170967
$ws.Cells.Item(5, 13).Value = 0.001310743449623976
$ws.Cells.Item(5, 14).Value = 0.001228004965801997
$ws.Cells.Item(5, 15).Value = 0.001228004965801997
$ws.Cells.Item(5, 16).Value = 0.001228004965801997
$ws.Cells.Item(5, 17).Value = 0.001199176269367962
$ws.Cells.Item(5, 18).Value = 0.00118768483859957
$ws.Cells.Item(5, 19).Value = 0.001179128012294062
$ws.Cells.Item(5, 20).Value = 0.001154071151656085
$ws.Cells.Item(5, 21).Value = 0.001135182204253374
$ws.Cells.Item(5, 22).Value = 0.001119167386653371
$ws.Cells.Item(5, 23).Value = 0.001099247544057896
$ws.Cells.Item(5, 24).Value = 0.001099247544057896
$ws.Cells.Item(5, 25).Value = 0.001089977932104191
$ws.Cells.Item(6, 3).Value = 0.3089587688446045
$ws.Cells.Item(6, 5).Value = 54.82122733192409
$ws.Cells.Item(6, 6).Value = 0.00193708617319336
$ws.Cells.Item(6, 7).Value = 0.001567339396814643
$ws.Cells.Item(6, 8).Value = 0.001472718318202052
$ws.Cells.Item(6, 9).Value = 0.001375652742961917
$ws.Cells.Item(6, 10).Value = 0.001360027046507985
$ws.Cells.Item(6, 11).Value = 0.001319507424841663
$ws.Cells.Item(6, 12).Value = 0.001319507424841663
$ws.Cells.Item(6, 13).Value = 0.001219204131700101
$ws.Cells.Item(6, 14).Value = 0.001219204131700101
$ws.Cells.Item(6, 15).Value = 0.001206774253131083
$ws.Cells.Item(6, 16).Value = 0.001176093531518318
$ws.Cells.Item(6, 17).Value = 0.001176093531518318
$ws.Cells.Item(6, 18).Value = 0.001143791534950704
$ws.Cells.Item(6, 19).Value = 0.001131665803806176
$ws.Cells.Item(6, 20).Value = 0.001120482994711816
$ws.Cells.Item(6, 21).Value = 0.001113034280203392
$ws.Cells.Item(6, 22).Value = 0.001094506547793344
$ws.Cells.Item(6, 23).Value = 0.00109305207247858
$ws.Cells.Item(6, 24).Value = 0.001088394536213544
$ws.Cells.Item(6, 25).Value = 0.001068639909004368
$ws.Cells.Item(7, 3).Value = 0.3329288959503174
$ws.Cells.Item(7, 5).Value = 55.45305847742566
$ws.Cells.Item(7, 6).Value = 0.00193708617319336
$ws.Cells.Item(7, 7).Value = 0.00160843379027147
$ws.Cells.Item(7, 8).Value = 0.001508323544057717
$ws.Cells.Item(7, 9).Value = 0.00147085894684559
$ws.Cells.Item(7, 10).Value = 0.001371826854793346
$ws.Cells.Item(7, 11).Value = 0.001353026232869847
$ws.Cells.Item(7, 12).Value = 0.00129613252471766
$ws.Cells.Item(7, 13).Value = 0.001163375587572966
$ws.Cells.Item(7, 14).Value = 0.001163375587572966
$ws.Cells.Item(7, 15).Value = 0.001163375587572966
$ws.Cells.Item(7, 16).Value = 0.001163375587572966
$ws.Cells.Item(7, 17).Value = 0.001163375587572966
$ws.Cells.Item(7, 18).Value = 0.001156715553782254
$ws.Cells.Item(7, 19).Value = 0.001148596602212179
$ws.Cells.Item(7, 20).Value = 0.001135982044601242
$ws.Cells.Item(7, 21).Value = 0.00111527730927445
$ws.Cells.Item(7, 22).Value = 0.001096330227938437
$ws.Cells.Item(7, 23).Value = 0.001096330227938437
$ws.Cells.Item(7, 24).Value = 0.001080956305602839
$ws.Cells.Item(7, 25).Value = 0.001080956305602839
$ws.Cells.Item(8, 3).Value = 0.2268500328063965
$ws.Cells.Item(8, 5).Value = 53.04076335701757
$ws.Cells.Item(8, 6).Value = 0.00193708617319336
$ws.Cells.Item(8, 7).Value = 0.0015608470781461
$ws.Cells.Item(8, 8).Value = 0.001420679607079712
$ws.Cells.Item(8, 9).Value = 0.001420679607079712
$ws.Cells.Item(8, 10).Value = 0.001355465243370035
$ws.Cells.Item(8, 11).Value = 0.001355465243370035
$ws.Cells.Item(8, 12).Value = 0.00123612566526615
$ws.Cells.Item(8, 13).Value = 0.00123612566526615
$ws.Cells.Item(8, 14).Value = 0.00123612566526615
$ws.Cells.Item(8, 15).Value = 0.00123612566526615
$ws.Cells.Item(8, 16).Value = 0.001175312739882922
$ws.Cells.Item(8, 17).Value = 0.00115670956468238
$ws.Cells.Item(8, 18).Value = 0.001142949423981591
$ws.Cells.Item(8, 19).Value = 0.001108828631174519
$ws.Cells.Item(8, 20).Value = 0.001096448671611447
$ws.Cells.Item(8, 21).Value = 0.001068596735772996
$ws.Cells.Item(8, 22).Value = 0.001053902763959515
$ws.Cells.Item(8, 23).Value = 0.001053902763959515
$ws.Cells.Item(8, 24).Value = 0.001041992099195964
$ws.Cells.Item(8, 25).Value = 0.001033933008908724
$ws.Cells.Item(9, 3).Value = 0.2073185443878174
$ws.Cells.Item(9, 5).Value = 53.75795821317661
$ws.Cells.Item(9, 6).Value = 0.00193708617319336
$ws.Cells.Item(9, 7).Value = 0.001583992595727725
$ws.Cells.Item(9, 8).Value = 0.00135923448286036
$ws.Cells.Item(9, 9).Value = 0.00135923448286036
$ws.Cells.Item(9, 10).Value = 0.001341167592764554
$ws.Cells.Item(9, 11).Value = 0.001341167592764554
$ws.Cells.Item(9, 12).Value = 0.001279973180365854
$ws.Cells.Item(9, 13).Value = 0.001185210842895596
$ws.Cells.Item(9, 14).Value = 0.001185210842895596
$ws.Cells.Item(9, 15).Value = 0.001161962389355659
$ws.Cells.Item(9, 16).Value = 0.001161962389355659
$ws.Cells.Item(9, 17).Value = 0.001159576680291361
$ws.Cells.Item(9, 18).Value = 0.001131076147370333
$ws.Cells.Item(9, 19).Value = 0.001106373517962029
$ws.Cells.Item(9, 20).Value = 0.001102305723478188
$ws.Cells.Item(9, 21).Value = 0.001097806337970324
$ws.Cells.Item(9, 22).Value = 0.001084765450392602
$ws.Cells.Item(9, 23).Value = 0.001063748249336683
$ws.Cells.Item(9, 24).Value = 0.00106007022577774
$ws.Cells.Item(9, 25).Value = 0.001047913415461532
$ws.Cells.Item(10, 3).Value = 0.239255428314209
$ws.Cells.Item(10, 5).Value = 54.71823375087115
$ws.Cells.Item(10, 6).Value = 0.001901045187600951
$ws.Cells.Item(10, 7).Value = 0.001716593631869141
$ws.Cells.Item(10, 8).Value = 0.001453525158727221
$ws.Cells.Item(10, 9).Value = 0.00140295595249291
$ws.Cells.Item(10, 10).Value = 0.001380974317021949
$ws.Cells.Item(10, 11).Value = 0.001380974317021949
$ws.Cells.Item(10, 12).Value = 0.00137031803002209
$ws.Cells.Item(10, 13).Value = 0.001244268751742688
$ws.Cells.Item(10, 14).Value = 0.001244268751742688
$ws.Cells.Item(10, 15).Value = 0.001229067786677871
$ws.Cells.Item(10, 16).Value = 0.001193838002500821
$ws.Cells.Item(10, 17).Value = 0.001193838002500821
$ws.Cells.Item(10, 18).Value = 0.001163407608780855
$ws.Cells.Item(10, 19).Value = 0.00115181574756884
$ws.Cells.Item(10, 20).Value = 0.001114666237604397
$ws.Cells.Item(10, 21).Value = 0.001103375675023001
$ws.Cells.Item(10, 22).Value = 0.001103375675023001
$ws.Cells.Item(10, 23).Value = 0.001087436560867582
$ws.Cells.Item(10, 24).Value = 0.00107539917713467
$ws.Cells.Item(10, 25).Value = 0.001066632236859086
$ws.Cells.Item(11, 3).Value = 0.2191896438598633
$ws.Cells.Item(11, 5).Value = 52.99324218674701
$ws.Cells.Item(11, 7).Value = 0.001625389986190719
$ws.Cells.Item(11, 8).Value = 0.001441146269323767
$ws.Cells.Item(11, 9).Value = 0.00136482833511409
$ws.Cells.Item(11, 10).Value = 0.001329357157207807
$ws.Cells.Item(11, 11).Value = 0.001329357157207807
$ws.Cells.Item(11, 12).Value = 0.001253597947444597
$ws.Cells.Item(11, 13).Value = 0.001253597947444597
$ws.Cells.Item(11, 14).Value = 0.001241764762831721
$ws.Cells.Item(11, 15).Value = 0.001241764762831721
$ws.Cells.Item(11, 16).Value = 0.001191641715421864
$ws.Cells.Item(11, 17).Value = 0.001146309840252627
$ws.Cells.Item(11, 18).Value = 0.001146309840252627
$ws.Cells.Item(11, 19).Value = 0.001110577520175511
$ws.Cells.Item(11, 20).Value = 0.001104657346617631
$ws.Cells.Item(11, 21).Value = 0.001100072085590861
$ws.Cells.Item(11, 22).Value = 0.001068409290372918
$ws.Cells.Item(11, 23).Value = 0.001053571629513561
$ws.Cells.Item(11, 24).Value = 0.001039166438432291
$ws.Cells.Item(11, 25).Value = 0.001033006670306959
